# DOMA-8525: ability to set Meter.isAutomatic via meters import
# Adds a new "Automatic" column (T) to the meter import example sheet,
# and normalizes the Unit / Account number / Meter number columns (B, D, F)
# plus the Verification date / Reading submission date columns (M, L)
# to plain text values instead of numbers / real dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 11

# --- 1. New "Automatic" header column (T) -----------------------------
# Give it the same look as the existing last header cell (S1) and then
# set its own text.
$ws.Cells.Item(1, 20).NumberFormat = $ws.Cells.Item(1, 19).NumberFormat
$ws.Cells.Item(1, 20).Value = "Automatic"

# Match the column width used by the neighbouring R:S columns.
$ws.Columns.Item(20).ColumnWidth = $ws.Columns.Item(19).ColumnWidth

# Data rows for column T stay empty, but still carry the row's normal
# (bordered/filled) style, mirroring column S on the same row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 20).NumberFormat = $ws.Cells.Item($r, 19).NumberFormat
}

# --- 2. Row data: Unit / Account number / Meter number become text ----
$unit = @{2="1";3="1";4="1";5="1";6="1";7="1";8="2";9="2";10="2";11="2"}
$account = @{2="111";3="111";4="111";5="111";6="111";7="111";8="222";9="222";10="222";11="222"}
$meterNumber = @{2="1";3="1";4="2";5="2";6="2";7="2";8="11";9="11";10="22";11="33"}

# Reading submission date (L) and Verification date (M) become plain text.
$readingDate = @{2="2021-12-20";3="2021-12-20";4="2021-12-20";5="2021-12-20";6="2021-11-20";7="2021-12-20";8="2021-12-21";9="2021-12-21";10="2021-12-20";11="2021-12-20"}
$verificationDate = @{2="2021-12-20";3="2021-12-20";4="2021-12-20";5="2021-12-20";6="2021-12-20";7="2021-12-20";8="2021-12-20";9="2021-12-20";10="2021-12-20";11="2021-12-20"}

for ($r = 2; $r -le $lastRow; $r++) {
    # Unit (B) / Account number (D) / Meter number (F): copy the text
    # formatting already used by their neighbour cells on the same row
    # (A/C/E), then assign the value as text.
    $ws.Cells.Item($r, 2).NumberFormat = $ws.Cells.Item($r, 1).NumberFormat
    $ws.Cells.Item($r, 2).Value = $unit[$r]

    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r, 3).NumberFormat
    $ws.Cells.Item($r, 4).Value = $account[$r]

    $ws.Cells.Item($r, 6).NumberFormat = $ws.Cells.Item($r, 5).NumberFormat
    $ws.Cells.Item($r, 6).Value = $meterNumber[$r]

    # Value 2 / Value 3 / Value 4 (I, J, K) keep their current (blank)
    # content but pick up the plain "General" numeric style used by
    # Value 1 (H) on the same row.
    $ws.Cells.Item($r, 9).NumberFormat = $ws.Cells.Item($r, 8).NumberFormat
    $ws.Cells.Item($r, 10).NumberFormat = $ws.Cells.Item($r, 8).NumberFormat
    $ws.Cells.Item($r, 11).NumberFormat = $ws.Cells.Item($r, 8).NumberFormat

    # Reading submission date (L): make sure it is textual, matching the
    # style already used by the Next verification date (N) column.
    $ws.Cells.Item($r, 12).NumberFormat = $ws.Cells.Item($r, 14).NumberFormat
    $ws.Cells.Item($r, 12).Value = $readingDate[$r]

    # Verification date (M): switch from a real date value to plain text.
    $ws.Cells.Item($r, 13).NumberFormat = $ws.Cells.Item($r, 14).NumberFormat
    $ws.Cells.Item($r, 13).Value = $verificationDate[$r]
}
